$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.769.19"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "2.313.45"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.91%  "
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.608"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.16"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.50%  "
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.37"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.98%  "
$ws.Range("E13").Value = "  -1.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.996"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.27%  "
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("D16").Value = "2.664.06"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").Value = "2.311.54"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").Value = "42.762.16"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("E19").Value = "  -1.14%  "
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -9.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.67%  "
$ws.Range("E23").Value = "  -1.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("E25").Value = "  +0.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +14.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.31"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.60"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.87%  "
$ws.Range("E31").Value = "  -1.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0874"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.65%  "
$ws.Range("E35").Value = "  -1.39%  "
$ws.Range("E36").Value = "  -0.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.65"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0357"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.66%  "
$ws.Range("E39").Value = "  +4.52%  "
$ws.Range("E40").Value = "  -2.08%  "
$ws.Range("E41").Value = "  -0.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "104.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "71.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.97%  "
$ws.Range("E44").Value = "  +1.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "112.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.12%  "
$ws.Range("D48").Value = "1.665.71"
$ws.Range("E48").Value = "  -1.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "76.86"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.65%  "
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("E51").Value = "  +0.50%  "
